$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H (old H -> I, old I -> J)
$ws.Columns("H").Insert()

# New header cell
$ws.Cells.Item(1,8).Value = "total_arqueo_ciego"

# Columns that must stay/become TEXT even though they look numeric/date-like:
# D (dates), H, I, J (decimal/operation counts with locale-specific formatting)
$textCols = 4,8,9,10
foreach ($col in $textCols) {
  $ws.Range($ws.Cells.Item(2,$col), $ws.Cells.Item(11,$col)).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2,2).Value = "V2"
$ws.Cells.Item(2,3).Value = "BAR"
$ws.Cells.Item(2,4).Value = "01/02/2025"
$ws.Cells.Item(2,5).Value = 8829
$ws.Cells.Item(2,7).Value = "EUROS"
$ws.Cells.Item(2,8).Value = "551.30"
$ws.Cells.Item(2,9).Value = "72,75"
$ws.Cells.Item(2,10).Value = "9"

# Row 3
$ws.Cells.Item(3,2).Value = "V2"
$ws.Cells.Item(3,3).Value = "BAR"
$ws.Cells.Item(3,4).Value = "01/02/2025"
$ws.Cells.Item(3,5).Value = 8829
$ws.Cells.Item(3,7).Value = "TARJETA VISA"
$ws.Cells.Item(3,8).Value = "891.60"
$ws.Cells.Item(3,9).Value = "202,10"
$ws.Cells.Item(3,10).Value = "24"

# Row 4
$ws.Cells.Item(4,2).Value = "V2"
$ws.Cells.Item(4,3).Value = "BAR"
$ws.Cells.Item(4,4).Value = "01/02/2025"
$ws.Cells.Item(4,5).Value = 8830
$ws.Cells.Item(4,7).Value = "EUROS"
$ws.Cells.Item(4,8).Value = "422.05"
$ws.Cells.Item(4,9).Value = "83,40"
$ws.Cells.Item(4,10).Value = "9"

# Row 5
$ws.Cells.Item(5,2).Value = "V2"
$ws.Cells.Item(5,3).Value = "BAR"
$ws.Cells.Item(5,4).Value = "01/02/2025"
$ws.Cells.Item(5,5).Value = 8830
$ws.Cells.Item(5,7).Value = "TARJETA VISA"
$ws.Cells.Item(5,8).Value = "867.00"
$ws.Cells.Item(5,9).Value = "139,80"
$ws.Cells.Item(5,10).Value = "14"

# Row 6
$ws.Cells.Item(6,2).Value = "V1"
$ws.Cells.Item(6,3).Value = "SERVIDOR TIENDA"
$ws.Cells.Item(6,4).Value = "01/02/2025"
$ws.Cells.Item(6,5).Value = 8828
$ws.Cells.Item(6,7).Value = "EUROS"
$ws.Cells.Item(6,8).Value = "1204.30"
$ws.Cells.Item(6,9).Value = "806,32"
$ws.Cells.Item(6,10).Value = "77"

# Row 7
$ws.Cells.Item(7,2).Value = "V1"
$ws.Cells.Item(7,3).Value = "SERVIDOR TIENDA"
$ws.Cells.Item(7,4).Value = "01/02/2025"
$ws.Cells.Item(7,5).Value = 8828
$ws.Cells.Item(7,7).Value = "TARJETA VISA"
$ws.Cells.Item(7,8).Value = "2231.23"
$ws.Cells.Item(7,9).Value = "2231,03"
$ws.Cells.Item(7,10).Value = "154"

# Row 8
$ws.Cells.Item(8,2).Value = "V1"
$ws.Cells.Item(8,3).Value = "SERVIDOR TIENDA"
$ws.Cells.Item(8,4).Value = "01/02/2025"
$ws.Cells.Item(8,5).Value = 8831
$ws.Cells.Item(8,7).Value = "EUROS"
$ws.Cells.Item(8,8).Value = "1130.48"
$ws.Cells.Item(8,9).Value = "702,76"
$ws.Cells.Item(8,10).Value = "79"

# Row 9
$ws.Cells.Item(9,2).Value = "V1"
$ws.Cells.Item(9,3).Value = "SERVIDOR TIENDA"
$ws.Cells.Item(9,4).Value = "01/02/2025"
$ws.Cells.Item(9,5).Value = 8831
$ws.Cells.Item(9,7).Value = "SMS"
$ws.Cells.Item(9,8).Value = "0.00"
$ws.Cells.Item(9,9).Value = "4,70"
$ws.Cells.Item(9,10).Value = "1"

# Row 10
$ws.Cells.Item(10,2).Value = "V1"
$ws.Cells.Item(10,3).Value = "SERVIDOR TIENDA"
$ws.Cells.Item(10,4).Value = "01/02/2025"
$ws.Cells.Item(10,5).Value = 8831
$ws.Cells.Item(10,7).Value = "TARJETA VISA"
$ws.Cells.Item(10,8).Value = "1801.60"
$ws.Cells.Item(10,9).Value = "1809,80"
$ws.Cells.Item(10,10).Value = "154"

# Row 11
$ws.Cells.Item(11,2).Value = "V2"
$ws.Cells.Item(11,3).Value = "BAR"
$ws.Cells.Item(11,4).Value = "01/02/2025"
$ws.Cells.Item(11,5).Value = 8833
$ws.Cells.Item(11,7).Value = "TARJETA VISA"
$ws.Cells.Item(11,8).Value = "1373.00"
$ws.Cells.Item(11,9).Value = "93,20"
$ws.Cells.Item(11,10).Value = "13"
